# Update the public EPEX Spot prices workbook with the newest day of data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add a new date column BE (09-aug) with 24 hourly
# values, mirroring the style of the preceding column (BD).
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header cell formatting (bold, bordered, centered) from BD1
# onto the new BE1 header cell before writing its value.
$wsSpot.Range("BD1").Copy()
$wsSpot.Range("BE1").PasteSpecial(-4122)
$wsSpot.Range("BE1").Value = "09-aug"

$spotValues = @{
    2  = 86.26000000000001
    3  = 70.13
    4  = 54.77
    5  = 42
    6  = 41.4
    7  = 44.03
    8  = 50.95
    9  = 47.81
    10 = 46.33
    11 = 15.34
    12 = -0.02
    13 = -1
    14 = -0.02
    15 = -4.74
    16 = -7.75
    17 = -3.96
    18 = -0.05
    19 = 3.52
    20 = 28.61
    21 = 69.98999999999999
    22 = 85.70999999999999
    23 = 78.95
    24 = 86.33
    25 = 78.91
}

foreach ($row in 2..25) {
    $wsSpot.Range("BE$row").Value = $spotValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append the next day's closing price as a new row.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date cell to be stored as literal text (matching the rest
# of column A) instead of letting Excel auto-convert it to a date
# serial number, then restore the default "Normal" style so no extra
# formatting is left behind.
$wsGaz.Range("A54").NumberFormat = "@"
$wsGaz.Range("A54").Value = "2025-08-07"
$wsGaz.Range("A54").Style = "Normal"
$wsGaz.Range("B54").Value = 32.175

# ---------------------------------------------------------------------
# Sheet "CO2": append the next day's closing price as a new row.
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A54").NumberFormat = "@"
$wsCO2.Range("A54").Value = "2025-08-07"
$wsCO2.Range("A54").Style = "Normal"
$wsCO2.Range("B54").Value = 71.15000000000001

Write-Output "done"
